$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# ------------------------------------------------------------------
# 1. Rebuild the hyperlinks: this runtime's per-item Hyperlink.Delete()
#    does not actually shrink the collection, so instead we clear every
#    hyperlink on the sheet and re-add only the six that should remain
#    (F2:F7), restoring the "Hyperlink" cell style they had before.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5438171")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5438740")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5438554")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5438369")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5438567")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5438092")
$ws.Range("F2:F7").Style = "Hyperlink"

# ------------------------------------------------------------------
# 2. Overwrite the data rows 2:7 with the new scrape results.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "2025-11-21 06:28:07"
$ws.Range("B2").Value = "【謝礼2,000円】AIに興味のあるエンジニアの方へ|45分だけお話を聞かせてください"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5438171"
$ws.Range("G2").Value = 295
$ws.Range("H2").Value = "🔥AI,Ai"

$ws.Range("A3").Value = "2025-11-21 06:28:07"
$ws.Range("B3").Value = "【データベース化】エクセル管理台帳の視覚化と検索機能強化"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5438740"
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = "◇管理"

$ws.Range("A4").Value = "2025-11-21 06:28:07"
$ws.Range("B4").Value = "【急募】東京の未来を形作る奉仕システム構築支援"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5438554"
$ws.Range("G4").Value = 40
$ws.Range("H4").ClearContents()

$ws.Range("A5").Value = "2025-11-21 06:28:07"
$ws.Range("B5").Value = "【教育分野】新プロジェクトのPM募集!企画整理とチーム構築"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5438369"
$ws.Range("G5").Value = 18
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = "2025-11-21 06:28:07"
$ws.Range("B6").Value = "【急募】instagramとSTORES連携で商品販売を実現したい"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5438567"
$ws.Range("G6").Value = 13
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = "2025-11-21 06:28:07"
$ws.Range("B7").Value = "PowerAutomate でWorepress記事を自動作成"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5438092"
$ws.Range("G7").Value = 13
$ws.Range("H7").ClearContents()

# ------------------------------------------------------------------
# 3. Remove the now-obsolete trailing rows (old rows 8:26).
# ------------------------------------------------------------------
$ws.Rows("8:26").Delete()

# ------------------------------------------------------------------
# 4. Resize columns B and H. Excel stores the OOXML "width" as the
#    COM ColumnWidth plus 5/6 (the default-font padding), so we back
#    that constant out to land exactly on the target stored widths
#    (B=45, H=12).
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 45 - 5/6
$ws.Columns("H").ColumnWidth = 12 - 5/6
